$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("buildings")

# Rename the "gas_consumption_kwhpa" header (column S) to
# "gas_consumption_m3pa" -> gas is now tracked in m3 per annum instead of kWh per annum.
$ws.Range("S1").Value = "gas_consumption_m3pa"
